# Replace the field `{ m:Sequence{...} }` (a complex field made of
# fldChar begin/instrText runs/fldChar end) in the second paragraph with
# plain literal text runs spelling out the same content (now wrapped in
# an extra pair of braces), keeping the existing `_GoBack` bookmark in
# place. This mirrors the parser switch to
# TokenIteratorFieldRewriterSplit, which emits template tokens as plain
# `w:t` runs instead of Word field codes.

$d = $word.ActiveDocument

# Locate the paragraph that holds the "m:Sequence{...}" field (the field
# code text is not visible through Range.Text, so a plain text Find won't
# see it -- instead find the paragraph whose Range contains a Field).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p.Range
    }
}

# Rebuild the paragraph as plain text runs (no field, no w:instrText),
# reproducing the exact run split from the target markup and keeping the
# `_GoBack` bookmark between the 6th and 7th run.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="006E189B" w:rsidRDefault="00B6779C">' +
    '<w:r><w:t xml:space="preserve">{m:Sequence{' + [char]39 + 'Some text in a first </w:t></w:r>' +
    '<w:r w:rsidR="006E6A44"><w:t>section of text</w:t></w:r>' +
    '<w:r><w:t>' + [char]39 + ', ' + [char]39 + 'new</w:t></w:r>' +
    '<w:r w:rsidR="006E6A44"><w:t>TextWrapping</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">' + [char]39 + '.asPagination(), ' + [char]39 + 'Some text in a second </w:t></w:r>' +
    '<w:r w:rsidR="006E6A44"><w:t>section of text</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve">' + [char]39 + '}}</w:t></w:r>' +
    '</w:p>'

$target.InsertXML($xml)
